$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values in columns D and E stay as text
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "291.07"
$ws.Range("E2").Value = "-8.12%"
$ws.Range("D3").Value = "40.42"
$ws.Range("E3").Value = "-1.57%"
$ws.Range("D4").Value = "5.013"
$ws.Range("E4").Value = "-2.68%"
$ws.Range("D5").Value = "0.07293"
$ws.Range("E5").Value = "-4.52%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "1.557"
$ws.Range("E6").Value = "-7.10%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "0.9220"
$ws.Range("E7").Value = "-1.23%"
$ws.Range("B8").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C8").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D8").Value = "0.1159"
$ws.Range("E8").Value = "-7.41%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "0.1733"
$ws.Range("E9").Value = "-6.06%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "0.08605"
$ws.Range("E10").Value = "-5.01%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "0.04172"
$ws.Range("E11").Value = "0.56%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "0.1052"
$ws.Range("E12").Value = "-0.38%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "0.001269"
$ws.Range("E13").Value = "-1.40%"
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").Value = "0.005803"
$ws.Range("E14").Value = "-2.55%"
$ws.Range("B15").Value = "LEO"
$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D15").Value = "3.398"
$ws.Range("E15").Value = "1.42%"
$ws.Range("B16").Value = "GateToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D16").Value = "4.281"
$ws.Range("E16").Value = "-0.84%"
$ws.Range("D18").Value = "0.3278"
$ws.Range("E18").Value = "-2.49%"
$ws.Range("D19").Value = "7.822"
$ws.Range("E19").Value = "-6.48%"
$ws.Range("E20").Value = "2.51%"
$ws.Range("D21").Value = "0.2882"
$ws.Range("E21").Value = "0.44%"
$ws.Range("E22").Value = "-4.48%"
$ws.Range("D23").Value = "0.001261"
$ws.Range("E23").Value = "-0.82%"
$ws.Range("D24").Value = "0.003798"
$ws.Range("E24").Value = "-6.94%"
$ws.Range("E25").Value = "0.44%"
$ws.Range("D26").Value = "0.0003724"
$ws.Range("D38").Value = "0.02307"
$ws.Range("E38").Value = "-7.71%"
$ws.Range("D39").Value = "0.04961"
$ws.Range("E39").Value = "-5.31%"
$ws.Range("E40").Value = "211.57%"
$ws.Range("D41").Value = "0.007695"
$ws.Range("E41").Value = "-0.77%"
$ws.Range("D42").Value = "0.1274"
$ws.Range("E42").Value = "-1.73%"
$ws.Range("D43").Value = "0.007370"
$ws.Range("E43").Value = "4.20%"
$ws.Range("D44").Value = "0.007085"
$ws.Range("E44").Value = "-14.00%"
$ws.Range("D45").Value = "0.3153"
$ws.Range("E45").Value = "-0.34%"
$ws.Range("D46").Value = "0.00006428"
$ws.Range("E46").Value = "-3.63%"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("D48").Value = "0.01633"
$ws.Range("E48").Value = "-95.15%"
$ws.Range("E49").Value = "-0.37%"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("D51").Value = "0.0002001"
